# Rename the worksheet from "Sheet1" to "test_cases"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "test_cases"

# Remove the bold/bordered/centered header styling that was applied to A1:E1
# (reverts the header row back to the default, unstyled cell format)
$ws.Range("A1:E1").ClearFormats()

# Update D5 so the contexts list is stored as a JSON-style array string
# instead of the old pipe-delimited string
$ws.Range("D5").Value = '["Test-Driven Development", "Write tests first", "Red-Green-Refactor cycle"]'
